$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.937.22'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.910.15'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '355.91'
$ws.Range("E5").Value = '  +0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.53'
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.566'
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.627'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.69'
$ws.Range("E10").Value = '  -3.24%  '
$ws.Range("E11").Value = '  +1.50%  '
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.46'
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.75'
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.367.52'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.902.80'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.979'
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.929.20'
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.41'
$ws.Range("E19").Value = '  +3.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.51'
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.88'
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0976'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.38'
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.72'
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.80'
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.182'
$ws.Range("E26").Value = '  +8.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.82'
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.60'
$ws.Range("E28").Value = '  +15.18%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +8.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.45'
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.34'
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.17'
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.21'
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '52.13'
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0442'
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.20'
$ws.Range("E38").Value = '  -3.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.13'
$ws.Range("E39").Value = '  -3.21%  '
$ws.Range("E40").Value = '  -4.13%  '
$ws.Range("E41").Value = '  -5.02%  '
$ws.Range("E42").Value = '  +2.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.92'
$ws.Range("E43").Value = '  -4.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.60'
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.17'
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("E46").Value = '  -6.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.44'
$ws.Range("E47").Value = '  -2.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.125.26'
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.251'
$ws.Range("E49").Value = '  -4.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0336'
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.913'
$ws.Range("E51").Value = '  -5.29%  '
